# Update the IBAN / account-holder list on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (IBAN, account holder name) that replaces the old rows 2-7
# and extends the table down to row 11.
$data = @(
    @("NL86RABO6333227641", "Sally Snozcumber"),
    @("NL06ABNA5558304928", "Roy Olsson"),
    @("NL36INGB2682297498", "Barry Grey"),
    @("NL10RABO9837080566", "Chloe Donaldson"),
    @("NL89INGB6034837898", "Alison Blackman"),
    @("NL57ABNA2454554658", "Gemma Parkes"),
    @("NL86INGB4110487447", "Suzanne Blast"),
    @("NL75ABNA9372718300", "Sally Lakeman"),
    @("NL23RABO5299017782", "Hannah Connor"),
    @("NL22ABNA5206019070", "Morwenna Zeus")
)

# Write the values for rows 2..11 (existing rows get overwritten, rows
# 8-11 are brand new).
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Make sure the newly created rows pick up the same formatting as the
# rest of the data rows.
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B11").PasteSpecial(-4122)

# New IBANs are one character longer than the old ones, so column A is a
# touch wider (best-fit) now.
$ws.Columns("A").ColumnWidth = 18.6667

# Mimic the selection left behind after pasting/typing the new values.
$ws.Range("A2:B11").Select() | Out-Null

# Page setup was touched as well.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
